# Sprint_1_Logs.xlsx - "Added to devinfo and help screens"
# Fills in the previously-empty "3/10/2019: 10pm" (column J) and
# "3/11/2019: 6pm" (column K) devlog answers for every question/person row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-IfNotNull($colLetter, $row, $text) {
    if ($null -ne $text) {
        $ws.Range("$colLetter$row").Value = $text
    }
}

# row -> [J text, K text (or $null if untouched)]
$rows = @{
    2  = @("Got intents working and linked all of the screens I have built thus far. Also renamed layout components in the format of Activity_NameType. This improved the readability in the Java code when mapping java objets to their corrects IDs.", $null)
    3  = @("I worked on the XML and Java file for the main menu, sign up page, and about page.", "I worked on the intents and Implementation for the main menu, sign up page and about page")
    4  = @("Created the layout and some functionality for the review page", "Finished the current layout for the help page, about page, and review page")
    5  = @("Continued attempting to port a Unity project to Android studio", "Successfully brought a Unity project over to Android studio and ran said project")
    6  = @("Finish planning all topics of discussion for our 'merge meeting', update the github, and continue to look into firebase", $null)
    7  = @("I will work on my assigned issues. Every Issue assigned to me as well as the design and layout of every button and image in all pages for all issues of all assignee.", "I will work on my assigned issues. Every issue assigned to me as well as the design and layout of every button and image in all pages for all issues of all assignee.")
    8  = @("Finish functionality of review, help, and info pages", "Continue to work on assigned issues")
    9  = @("Continue looking into porting a Unity project to Android Studio", "Mapping character movement to UI buttons and integrating the Unity project into a pre-existing Android Studio project")
    10 = @("Not currently", $null)
    11 = @("No, I am still working with all cylinders pumping", "No, I am still working with all cylinders pumping")
    12 = @("Travel ", "No")
    13 = @("Several errors have prevented me from successfully running a Unity proect in Android Studio", "Nothing is currently getting in the way of my work")
    14 = @("The use of intentions", $null)
    15 = @("Working together is better than alone!", "Working together is usually better than alone!")
    16 = @("Learning different ways of adjusting EditText and Ratings", "How to use and alter the rating bar and clear them after submission is made")
    17 = @("Learned more about how to run a Unity project in Andriod Studio", "Learned how to port a Unity project into Android Studio")
    18 = @("Changing the format of layout widget names will need to happen at some point to standardize them", $null)
    19 = @("A consistent color scheme/theme will be a necessary addition. However, most likely will not occur until the end of the project", "No")
    20 = @("Not yet", "No")
    21 = @("No changes currently have to be made to the current plan for the project", "No changes currently have to be made to the current plan for the project")
}

foreach ($r in $rows.Keys | Sort-Object) {
    $vals = $rows[$r]
    Set-IfNotNull "J" $r $vals[0]
    Set-IfNotNull "K" $r $vals[1]
}

# Row 7 grew taller once the new J/K answers wrapped onto more lines.
$ws.Rows(7).RowHeight = 88

# Leave the selection where the author ended up after typing the last entry.
$ws.Range("L14").Select()
